$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "PRODUCTO" column header in K1
$ws.Range("K1").Value = "PRODUCTO"

# Fill K2:K239 with "GIRASOL" for every data row
$ws.Range("K2:K239").Value = "GIRASOL"
